$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers) ---
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Type"
$ws.Range("C1").Value = "Base_Prob"
$ws.Range("D1").Value = "Min_Prob"
$ws.Range("E1").Value = "Max_Prob"
$ws.Range("F1").Value = "Scaling"
$ws.Range("G1").Value = "Impacts"

# --- Row 2 (earthquake) ---
$ws.Range("A2").Value = "earthquake"
$ws.Range("B2").Value = "deterministic"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0.0001
$ws.Range("E2").Value = 0.05
$ws.Range("G2").Value = "[population -10%, landArea -5%]"

# --- New Impacts/Scaling text values (entered in this order) ---
$ws.Range("F4").Value = "[farmLand 1%]"
$ws.Range("F2").Value = "[landArea 0.5%]"
$ws.Range("A3").Value = "hurricane"
$ws.Range("F3").Value = "[landArea 0.15%]"
$ws.Range("G3").Value = "[population -5%, timber -20%, housing -25%]"
$ws.Range("A4").Value = "drought"

# --- Row 3 (hurricane, was worker_strike) ---
$ws.Range("B3").Value = "natural"
$ws.Range("C3").Value = 0.02
$ws.Range("D3").Value = 0.0005
$ws.Range("D3").NumberFormat = "0.00%"
$ws.Range("E3").Value = 0.065
$ws.Range("E3").NumberFormat = "0.00%"

# --- Row 4 (drought, was drought_natural) ---
$ws.Range("B4").Value = "natural"
$ws.Range("C4").Value = 0.03
$ws.Range("D4").Value = 0.02
$ws.Range("D4").NumberFormat = "0%"
$ws.Range("E4").Value = 0.1
$ws.Range("E4").NumberFormat = "0%"
$ws.Range("G4").Value = "[farmLand -25%, landArea -5%]"

# --- Row 5 (was drought_induced) - fully cleared, no leftover style ---
$ws.Range("A5:G5").ClearContents()
$ws.Range("A5:G5").Style = "Normal"

# --- Column G width widened to fit the longer new impact text ---
$ws.Columns.Item(7).ColumnWidth = 38

# --- Re-anchor selection away from the row that was just emptied ---
$ws.Range("A1").Select()
